# Change point detection algorithm
# Strip the "_upd" suffix from the sample identifiers in column A,
# normalizing "FInf" to "Finf" where needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "C23F07"
$ws.Range("A3").Value = "C24F07"
$ws.Range("A4").Value = "C25F07"
$ws.Range("A5").Value = "C26Finf"
$ws.Range("A6").Value = "C27Finf"
$ws.Range("A7").Value = "C28Finf"
$ws.Range("A8").Value = "C29Finf"
$ws.Range("A9").Value = "C30Finf"
$ws.Range("A10").Value = "C31Finf"
$ws.Range("A11").Value = "C32Finf"
$ws.Range("A12").Value = "C33Finf"
$ws.Range("A13").Value = "C34Finf"
